$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> new value for column F (dSF)
$updates = @{
    11 = 0
    15 = 1
    16 = 3
    20 = 1
    27 = 0
    28 = -1
    33 = -1
    39 = 2
    40 = 0
    50 = 0
    52 = 0
    53 = 2
    64 = 0
    68 = 0
    70 = -6
    71 = -2
    78 = -5
    82 = -12
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
